# Update the "Bonus" column (K) values for rows 2-15 with new random integer values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 258
    3  = 275
    4  = 242
    5  = 138
    6  = 115
    7  = 249
    8  = 154
    9  = 144
    10 = 179
    11 = 115
    12 = 186
    13 = 186
    14 = 132
    15 = 214
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}
